# Update the cached "datetimeFigureOut" date fields (slide master + the
# "Première page" slide layout) from 2023-06-20 to 2023-06-21.
$p = $ppt.ActivePresentation

$master = $p.SlideMaster
$masterDateShape = $master.Shapes.Item(3)
$masterDateShape.TextFrame.TextRange.Text = "2023-06-21"

$layout = $master.CustomLayouts.Item(1)
$layoutDateShape = $layout.Shapes.Item(2)
$layoutDateShape.TextFrame.TextRange.Text = "2023-06-21"

# Update the "banner" placeholder text to the literal token "<banner>" on
# both slides that contain it.
$slide1 = $p.Slides.Item(1)
$slide1.Shapes.Item(6).TextFrame.TextRange.Text = "<banner>"

$slide2 = $p.Slides.Item(2)
$slide2.Shapes.Item(1).TextFrame.TextRange.Text = "<banner>"
